$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the columns that held the old "F0000xx" / name / description columns
# (old B..G) -- this shifts old H,I,J into the new B,C,D slots and keeps
# their original widths/styles intact.
$ws.Columns("B:G").Delete()

# Remove the old sample rows 3-9 (F000006 ... F000029), leaving only the header row.
$ws.Range("A3:A9").EntireRow.Delete()

# Re-label the header row for the new "parceiro / nc" relationship sheet.
# Values are written B,C,D first, then A, so the shared-string table ends up
# ordered the same way the source workbook has it.
$ws.Range("B1").Value = "nome_responsavel"
$ws.Range("C1").Value = "email_responsavel"
$ws.Range("D1").Value = "telefone_responsavel"
$ws.Range("A1").Value = "nome"

# Move the selection the way the saved workbook shows it.
[void]$ws.Range("C7").Select()

# Re-apply the autofilter over the new (narrower) data range.
[void]$ws.Range("A1:D418").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new range.
foreach ($n in $wb.Names) {
  if ($n.Name -eq "Planilha1!_FilterDatabase") {
    $n.RefersTo = "=Planilha1!`$A`$1:`$D`$418"
  }
}
